# Update the raw unit_price data in column D; column E (amount) recalculates
# via its existing qty*unit_price formula.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1500
$ws.Range("D3").Value = 2000
$ws.Range("D4").Value = 120
$ws.Range("D5").Value = 400

# Rows 3-5 didn't have the amount formula applied yet (they held stale
# hard-coded values) -- bring them in line with row 2's C*D formula.
$ws.Range("E3:E5").Formula = "=C3*D3"

# Reflect the edited range in the sheet's selection.
$ws.Range("E3:E5").Select()
